$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $old = $cell.Value2
    $new = $old -replace '^sequence/run_0673_samples/', ''
    $cell.Value2 = $new
}

$ws.Range("F16").Select()
